# Restore/update the value of C10 on the active sheet from 18 to 1,
# matching the revision restored in the commit ("Restored from revision ...").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
